$d = $word.ActiveDocument

$d.Content.Find.Execute("22-7=", $true, $true, $false, $false, $false, $true, 1, $false, "11+70=", 2) | Out-Null
$d.Content.Find.Execute("31-25=", $true, $true, $false, $false, $false, $true, 1, $false, "47-8=", 2) | Out-Null
$d.Content.Find.Execute("51+4=", $true, $true, $false, $false, $false, $true, 1, $false, "71-24=", 2) | Out-Null
$d.Content.Find.Execute("75-64=", $true, $true, $false, $false, $false, $true, 1, $false, "84-57=", 2) | Out-Null
$d.Content.Find.Execute("52-49=", $true, $true, $false, $false, $false, $true, 1, $false, "87-76=", 2) | Out-Null
$d.Content.Find.Execute("96-60=", $true, $true, $false, $false, $false, $true, 1, $false, "45-2=", 2) | Out-Null
$d.Content.Find.Execute("61-3=", $true, $true, $false, $false, $false, $true, 1, $false, "69-46=", 2) | Out-Null
$d.Content.Find.Execute("46+28=", $true, $true, $false, $false, $false, $true, 1, $false, "61-13=", 2) | Out-Null
$d.Content.Find.Execute("8+63=", $true, $true, $false, $false, $false, $true, 1, $false, "65-31=", 2) | Out-Null
$d.Content.Find.Execute("2+88=", $true, $true, $false, $false, $false, $true, 1, $false, "49+6=", 2) | Out-Null
$d.Content.Find.Execute("20+20=", $true, $true, $false, $false, $false, $true, 1, $false, "72-62=", 2) | Out-Null
$d.Content.Find.Execute("68-27=", $true, $true, $false, $false, $false, $true, 1, $false, "40+42=", 2) | Out-Null
$d.Content.Find.Execute("76-20=", $true, $true, $false, $false, $false, $true, 1, $false, "58+3=", 2) | Out-Null
$d.Content.Find.Execute("75-70=", $true, $true, $false, $false, $false, $true, 1, $false, "87+8=", 2) | Out-Null
$d.Content.Find.Execute("60+25=", $true, $true, $false, $false, $false, $true, 1, $false, "42+16=", 2) | Out-Null
$d.Content.Find.Execute("30+40=", $true, $true, $false, $false, $false, $true, 1, $false, "13+20=", 2) | Out-Null
$d.Content.Find.Execute("10+80=", $true, $true, $false, $false, $false, $true, 1, $false, "48+2=", 2) | Out-Null
$d.Content.Find.Execute("80-54=", $true, $true, $false, $false, $false, $true, 1, $false, "16-4=", 2) | Out-Null
$d.Content.Find.Execute("35-33=", $true, $true, $false, $false, $false, $true, 1, $false, "26+28=", 2) | Out-Null
$d.Content.Find.Execute("40+33=", $true, $true, $false, $false, $false, $true, 1, $false, "86-72=", 2) | Out-Null
$d.Content.Find.Execute("39-20=", $true, $true, $false, $false, $false, $true, 1, $false, "51-43=", 2) | Out-Null
$d.Content.Find.Execute("35+27=", $true, $true, $false, $false, $false, $true, 1, $false, "76+13=", 2) | Out-Null
$d.Content.Find.Execute("25-18=", $true, $true, $false, $false, $false, $true, 1, $false, "44-6=", 2) | Out-Null
$d.Content.Find.Execute("8+35=", $true, $true, $false, $false, $false, $true, 1, $false, "97-37=", 2) | Out-Null
$d.Content.Find.Execute("91-63=", $true, $true, $false, $false, $false, $true, 1, $false, "40-7=", 2) | Out-Null
$d.Content.Find.Execute("82-46=", $true, $true, $false, $false, $false, $true, 1, $false, "65-26=", 2) | Out-Null
$d.Content.Find.Execute("24+24=", $true, $true, $false, $false, $false, $true, 1, $false, "89-14=", 2) | Out-Null
$d.Content.Find.Execute("60+3=", $true, $true, $false, $false, $false, $true, 1, $false, "54-6=", 2) | Out-Null
$d.Content.Find.Execute("93-75=", $true, $true, $false, $false, $false, $true, 1, $false, "22+26=", 2) | Out-Null
$d.Content.Find.Execute("44-32=", $true, $true, $false, $false, $false, $true, 1, $false, "84-56=", 2) | Out-Null
$d.Content.Find.Execute("7+36=", $true, $true, $false, $false, $false, $true, 1, $false, "10+74=", 2) | Out-Null
$d.Content.Find.Execute("76-41=", $true, $true, $false, $false, $false, $true, 1, $false, "89-56=", 2) | Out-Null
$d.Content.Find.Execute("22-11=", $true, $true, $false, $false, $false, $true, 1, $false, "43-31=", 2) | Out-Null
$d.Content.Find.Execute("96-49=", $true, $true, $false, $false, $false, $true, 1, $false, "31+58=", 2) | Out-Null
$d.Content.Find.Execute("77+13=", $true, $true, $false, $false, $false, $true, 1, $false, "23+59=", 2) | Out-Null
$d.Content.Find.Execute("45+20=", $true, $true, $false, $false, $false, $true, 1, $false, "43+54=", 2) | Out-Null
$d.Content.Find.Execute("85-67=", $true, $true, $false, $false, $false, $true, 1, $false, "44+12=", 2) | Out-Null
$d.Content.Find.Execute("35+53=", $true, $true, $false, $false, $false, $true, 1, $false, "62-55=", 2) | Out-Null
$d.Content.Find.Execute("56+3=", $true, $true, $false, $false, $false, $true, 1, $false, "18+22=", 2) | Out-Null
$d.Content.Find.Execute("68-44=", $true, $true, $false, $false, $false, $true, 1, $false, "60+11=", 2) | Out-Null
$d.Content.Find.Execute("45-13=", $true, $true, $false, $false, $false, $true, 1, $false, "79-56=", 2) | Out-Null
$d.Content.Find.Execute("77-5=", $true, $true, $false, $false, $false, $true, 1, $false, "99-27=", 2) | Out-Null
$d.Content.Find.Execute("66-64=", $true, $true, $false, $false, $false, $true, 1, $false, "94-11=", 2) | Out-Null
$d.Content.Find.Execute("16+18=", $true, $true, $false, $false, $false, $true, 1, $false, "78+9=", 2) | Out-Null
$d.Content.Find.Execute("86+4=", $true, $true, $false, $false, $false, $true, 1, $false, "73-54=", 2) | Out-Null
$d.Content.Find.Execute("47+44=", $true, $true, $false, $false, $false, $true, 1, $false, "3+8=", 2) | Out-Null
$d.Content.Find.Execute("2+55=", $true, $true, $false, $false, $false, $true, 1, $false, "87-53=", 2) | Out-Null
$d.Content.Find.Execute("51+36=", $true, $true, $false, $false, $false, $true, 1, $false, "3+28=", 2) | Out-Null
$d.Content.Find.Execute("21+74=", $true, $true, $false, $false, $false, $true, 1, $false, "35-2=", 2) | Out-Null
$d.Content.Find.Execute("27+24=", $true, $true, $false, $false, $false, $true, 1, $false, "98-11=", 2) | Out-Null
$d.Content.Find.Execute("71-54=", $true, $true, $false, $false, $false, $true, 1, $false, "57+3=", 2) | Out-Null
$d.Content.Find.Execute("79-68=", $true, $true, $false, $false, $false, $true, 1, $false, "59-53=", 2) | Out-Null
$d.Content.Find.Execute("50-27=", $true, $true, $false, $false, $false, $true, 1, $false, "99-68=", 2) | Out-Null
$d.Content.Find.Execute("21+51=", $true, $true, $false, $false, $false, $true, 1, $false, "58-5=", 2) | Out-Null
$d.Content.Find.Execute("67+20=", $true, $true, $false, $false, $false, $true, 1, $false, "80-32=", 2) | Out-Null
$d.Content.Find.Execute("99-18=", $true, $true, $false, $false, $false, $true, 1, $false, "91-87=", 2) | Out-Null
$d.Content.Find.Execute("64-1=", $true, $true, $false, $false, $false, $true, 1, $false, "4+26=", 2) | Out-Null
$d.Content.Find.Execute("7+67=", $true, $true, $false, $false, $false, $true, 1, $false, "3+11=", 2) | Out-Null
$d.Content.Find.Execute("74+12=", $true, $true, $false, $false, $false, $true, 1, $false, "25+12=", 2) | Out-Null
$d.Content.Find.Execute("23+24=", $true, $true, $false, $false, $false, $true, 1, $false, "64+12=", 2) | Out-Null
$d.Content.Find.Execute("38+3=", $true, $true, $false, $false, $false, $true, 1, $false, "49-28=", 2) | Out-Null
$d.Content.Find.Execute("99-83=", $true, $true, $false, $false, $false, $true, 1, $false, "58-10=", 2) | Out-Null
$d.Content.Find.Execute("51-29=", $true, $true, $false, $false, $false, $true, 1, $false, "85-2=", 2) | Out-Null
$d.Content.Find.Execute("8+85=", $true, $true, $false, $false, $false, $true, 1, $false, "20+5=", 2) | Out-Null
$d.Content.Find.Execute("3+7=", $true, $true, $false, $false, $false, $true, 1, $false, "30+6=", 2) | Out-Null
$d.Content.Find.Execute("21+62=", $true, $true, $false, $false, $false, $true, 1, $false, "63+30=", 2) | Out-Null
$d.Content.Find.Execute("71-8=", $true, $true, $false, $false, $false, $true, 1, $false, "55+34=", 2) | Out-Null
$d.Content.Find.Execute("57-54=", $true, $true, $false, $false, $false, $true, 1, $false, "23+30=", 2) | Out-Null
$d.Content.Find.Execute("70+1=", $true, $true, $false, $false, $false, $true, 1, $false, "75-17=", 2) | Out-Null
$d.Content.Find.Execute("54-49=", $true, $true, $false, $false, $false, $true, 1, $false, "47-30=", 2) | Out-Null
$d.Content.Find.Execute("78-56=", $true, $true, $false, $false, $false, $true, 1, $false, "33+28=", 2) | Out-Null
$d.Content.Find.Execute("28+40=", $true, $true, $false, $false, $false, $true, 1, $false, "72+13=", 2) | Out-Null
$d.Content.Find.Execute("5+85=", $true, $true, $false, $false, $false, $true, 1, $false, "31+40=", 2) | Out-Null
$d.Content.Find.Execute("94-89=", $true, $true, $false, $false, $false, $true, 1, $false, "7+58=", 2) | Out-Null
$d.Content.Find.Execute("71+17=", $true, $true, $false, $false, $false, $true, 1, $false, "85+2=", 2) | Out-Null
$d.Content.Find.Execute("78-45=", $true, $true, $false, $false, $false, $true, 1, $false, "77-18=", 2) | Out-Null
$d.Content.Find.Execute("37+30=", $true, $true, $false, $false, $false, $true, 1, $false, "37-9=", 2) | Out-Null
$d.Content.Find.Execute("34+20=", $true, $true, $false, $false, $false, $true, 1, $false, "99-92=", 2) | Out-Null
$d.Content.Find.Execute("40+32=", $true, $true, $false, $false, $false, $true, 1, $false, "84-25=", 2) | Out-Null
$d.Content.Find.Execute("50-20=", $true, $true, $false, $false, $false, $true, 1, $false, "77+10=", 2) | Out-Null
$d.Content.Find.Execute("19-6=", $true, $true, $false, $false, $false, $true, 1, $false, "64+18=", 2) | Out-Null
$d.Content.Find.Execute("13+78=", $true, $true, $false, $false, $false, $true, 1, $false, "50+19=", 2) | Out-Null
$d.Content.Find.Execute("29+7=", $true, $true, $false, $false, $false, $true, 1, $false, "38+2=", 2) | Out-Null
$d.Content.Find.Execute("36-26=", $true, $true, $false, $false, $false, $true, 1, $false, "96-40=", 2) | Out-Null
$d.Content.Find.Execute("79-0=", $true, $true, $false, $false, $false, $true, 1, $false, "79-59=", 2) | Out-Null
$d.Content.Find.Execute("38-9=", $true, $true, $false, $false, $false, $true, 1, $false, "81-30=", 2) | Out-Null
$d.Content.Find.Execute("21-14=", $true, $true, $false, $false, $false, $true, 1, $false, "81-48=", 2) | Out-Null
$d.Content.Find.Execute("68+9=", $true, $true, $false, $false, $false, $true, 1, $false, "84-2=", 2) | Out-Null
$d.Content.Find.Execute("93-7=", $true, $true, $false, $false, $false, $true, 1, $false, "27-16=", 2) | Out-Null
$d.Content.Find.Execute("10-8=", $true, $true, $false, $false, $false, $true, 1, $false, "43+1=", 2) | Out-Null
$d.Content.Find.Execute("86-50=", $true, $true, $false, $false, $false, $true, 1, $false, "48+34=", 2) | Out-Null
$d.Content.Find.Execute("55+23=", $true, $true, $false, $false, $false, $true, 1, $false, "8+49=", 2) | Out-Null
$d.Content.Find.Execute("8+12=", $true, $true, $false, $false, $false, $true, 1, $false, "35+29=", 2) | Out-Null
$d.Content.Find.Execute("50-0=", $true, $true, $false, $false, $false, $true, 1, $false, "17+49=", 2) | Out-Null
$d.Content.Find.Execute("81+4=", $true, $true, $false, $false, $false, $true, 1, $false, "69-55=", 2) | Out-Null
$d.Content.Find.Execute("45+5=", $true, $true, $false, $false, $false, $true, 1, $false, "23+61=", 2) | Out-Null
$d.Content.Find.Execute("71+20=", $true, $true, $false, $false, $false, $true, 1, $false, "77+14=", 2) | Out-Null
$d.Content.Find.Execute("26+26=", $true, $true, $false, $false, $false, $true, 1, $false, "38+14=", 2) | Out-Null
$d.Content.Find.Execute("92-12=", $true, $true, $false, $false, $false, $true, 1, $false, "95-49=", 2) | Out-Null
$d.Content.Find.Execute("56-1=", $true, $true, $false, $false, $false, $true, 1, $false, "34+55=", 2) | Out-Null
